# CDUN - (Horario) Modificar bloque horario
# "Terminé el narrativo" - fill in the use-case narrative template with
# the actual content for the "Modificar Bloque Horario" use case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Referencias
$ws.Range("F5").Value = "Administrar Horario"

# Precondición
$ws.Range("E6").Value = "El bloque horario ya debe existir"

# Postcondición
$ws.Range("E7").Value = "Bloque horario modificado correctamente en el sistema."

# Autor / fecha / version
$ws.Range("E8").Value = "Iván Zamorano"
$ws.Range("G8").NumberFormat = "mm-dd-yy"
$ws.Range("G8").Value = (Get-Date -Year 2015 -Month 4 -Day 17 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "1.0"

# Propósito
$ws.Range("B11").Value = "Modificar un bloque horario que se encuentra registrado en el sistema."

# Resumen
$ws.Range("B14").Value = "Modifica las opciones de un bloque horario existente en el sistema, ya sea por un error o una modificación arbitraria"

# Curso normal (Básico)
$ws.Range("C17").Value = "Actor Admin: Se selecciona el bloque horario a modificar"
$ws.Range("G17").Value = "Actor Admin: Se seleccionan los cambios a modificar"
$ws.Range("C18").Value = "Actor Admin: Se guardan los cambios"
$ws.Range("G18").Value = "El sistema valida los datos y los guarda"

# Cursos alternos
$ws.Range("B21").Value = "3.a"
$ws.Range("C21").Value = "Los cambios no se pueden efectuar debido a que se ingresó un campo no válido. Se pide reingreso del dato"

# Otros datos: Frecuencia esperada / Rendimiento, Importancia / Urgencia, Estado / Estabilidad
$ws.Range("E25").Value = "Media"
$ws.Range("E26").Value = "Media"
$ws.Range("H26").Value = "Media"
$ws.Range("E27").Value = "Sin implementar"
$ws.Range("H27").Value = "Alta"

# Selection moved to J27 (cosmetic, matches the saved cursor position)
$ws.Range("J27").Select()

# Page setup Dpi bump
$ws.PageSetup.Orientation = 1
